$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the lecture plan cells for week 15 (row 14) and week 16 (row 15):
# - C14: the 07.04 seminar is cancelled, replaced by a note about oral exams
# - C15: the 14.04 entry becomes "Seminar 3" (moved from 07.04)
# - D15: the 16.04 entry becomes "Practical information about exam" (moved from 14.04),
#        and the old 16.04 Oracle session text is dropped (Oracle session stays on D14 = 09.04)
$ws.Range("D15").Value = "16.04: <strong>Practical information about exam</strong> in Aud C"
$ws.Range("C15").Value = "14.04: <strong>Seminar 3</strong> in Aud C"
$ws.Range("C14").Value = "07.04:  <strong>No lecture (oral exams in MAB1)</strong>"

# Update the view state to match the saved selection/scroll position
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("D16").Select()
